# Applies scheduled-runner data refresh to Mandragora Profits workbook
# Updates computed profit columns (H:N) across ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4494
$ws.Range("I69").Value = 4910
$ws.Range("J69").Value = 4216.6665
$ws.Range("K69").Value = 14730
$ws.Range("L69").Value = 12649.9995
$ws.Range("M69").Value = -13856
$ws.Range("N69").Value = -14397.9995

$ws.Range("H72").Value = 4494
$ws.Range("I72").Value = 4910
$ws.Range("J72").Value = 4216.6665
$ws.Range("K72").Value = 44190
$ws.Range("L72").Value = 37949.9985
$ws.Range("M72").Value = -39822
$ws.Range("N72").Value = -46685.9985

$ws.Range("H103").Value = 1022.26666
$ws.Range("I103").Value = 1188.75
$ws.Range("J103").Value = 832
$ws.Range("K103").Value = 3566.25
$ws.Range("L103").Value = 2496
$ws.Range("M103").Value = -2980.25
$ws.Range("N103").Value = -3668

$ws.Range("H133").Value = 42980
$ws.Range("J133").Value = 42980
$ws.Range("L133").Value = 42980
$ws.Range("N133").Value = -53100

$ws.Range("H137").Value = 2197.6858
$ws.Range("I137").Value = 3092.4666
$ws.Range("K137").Value = 9277.399800000001
$ws.Range("M137").Value = -6727.399800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1544.3889
$ws.Range("I74").Value = 1447.12
$ws.Range("J74").Value = 1765.4546
$ws.Range("K74").Value = 1447.12
$ws.Range("L74").Value = 1765.4546
$ws.Range("M74").Value = -573.1199999999999
$ws.Range("N74").Value = -3513.4546

$ws.Range("H77").Value = 1544.3889
$ws.Range("I77").Value = 1447.12
$ws.Range("J77").Value = 1765.4546
$ws.Range("K77").Value = 7235.599999999999
$ws.Range("L77").Value = 8827.273000000001
$ws.Range("M77").Value = -2867.599999999999
$ws.Range("N77").Value = -17563.273

$ws.Range("H122").Value = 2685
$ws.Range("I122").Value = 1599.1765
$ws.Range("J122").Value = 8838
$ws.Range("K122").Value = 4797.529500000001
$ws.Range("L122").Value = 26514
$ws.Range("M122").Value = -2347.529500000001
$ws.Range("N122").Value = -31414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8868.166999999999
$ws.Range("I7").Value = 12533.5
$ws.Range("J7").Value = 1537.5
$ws.Range("K7").Value = 12533.5
$ws.Range("L7").Value = 1537.5
$ws.Range("M7").Value = -12420.5
$ws.Range("N7").Value = -1763.5

$ws.Range("H20").Value = 43332.668
$ws.Range("J20").Value = 43332.668
$ws.Range("L20").Value = 43332.668
$ws.Range("N20").Value = -43804.668

$ws.Range("H30").Value = 43332.668
$ws.Range("J30").Value = 43332.668
$ws.Range("L30").Value = 43332.668
$ws.Range("N30").Value = -43514.668

$ws.Range("H31").Value = 2233.673
$ws.Range("I31").Value = 1493.0834
$ws.Range("J31").Value = 3900
$ws.Range("K31").Value = 1493.0834
$ws.Range("L31").Value = 3900
$ws.Range("M31").Value = -1198.0834
$ws.Range("N31").Value = -4490

$ws.Range("H34").Value = 2233.673
$ws.Range("I34").Value = 1493.0834
$ws.Range("J34").Value = 3900
$ws.Range("K34").Value = 1493.0834
$ws.Range("L34").Value = 3900
$ws.Range("M34").Value = -1291.0834
$ws.Range("N34").Value = -4304

$ws.Range("H58").Value = 1356348.9
$ws.Range("I58").Value = 1980.1818
$ws.Range("J58").Value = 3342756.2
$ws.Range("K58").Value = 1980.1818
$ws.Range("L58").Value = 3342756.2
$ws.Range("M58").Value = -1777.1818
$ws.Range("N58").Value = -3343162.2

$ws.Range("H128").Value = 43332.668
$ws.Range("J128").Value = 43332.668
$ws.Range("L128").Value = 43332.668
$ws.Range("N128").Value = -53292.668

$ws.Range("H136").Value = 1356348.9
$ws.Range("I136").Value = 1980.1818
$ws.Range("J136").Value = 3342756.2
$ws.Range("K136").Value = 5940.5454
$ws.Range("L136").Value = 10028268.6
$ws.Range("M136").Value = -3390.5454
$ws.Range("N136").Value = -10033368.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3377.6667
$ws.Range("J100").Value = 3377.6667
$ws.Range("L100").Value = 10133.0001
$ws.Range("N100").Value = -11755.0001

$ws.Range("H134").Value = 3276.3416
$ws.Range("I134").Value = 1523.5714
$ws.Range("J134").Value = 4185.185
$ws.Range("K134").Value = 4570.7142
$ws.Range("L134").Value = 12555.555
$ws.Range("M134").Value = 499.2857999999997
$ws.Range("N134").Value = -22695.555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H63").Value = 13229.117
$ws.Range("J63").Value = 13229.117
$ws.Range("L63").Value = 13229.117
$ws.Range("N63").Value = -14601.117

$ws.Range("H66").Value = 13229.117
$ws.Range("J66").Value = 13229.117
$ws.Range("L66").Value = 39687.351
$ws.Range("N66").Value = -46551.351

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H98").Value = 19950
$ws.Range("J98").Value = 19950
$ws.Range("L98").Value = 19950
$ws.Range("N98").Value = -25940

$ws.Range("H132").Value = 4169702.5
$ws.Range("I132").Value = 6947179.5
$ws.Range("J132").Value = 3487
$ws.Range("K132").Value = 20841538.5
$ws.Range("L132").Value = 10461
$ws.Range("M132").Value = -20839008.5
$ws.Range("N132").Value = -15521

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 971
$ws.Range("I9").Value = 106.5
$ws.Range("J9").Value = 2700
$ws.Range("K9").Value = 106.5
$ws.Range("L9").Value = 2700
$ws.Range("M9").Value = 117.5
$ws.Range("N9").Value = -3148

$ws.Range("H36").Value = 40000
$ws.Range("J36").Value = 40000
$ws.Range("L36").Value = 40000
$ws.Range("N36").Value = -41124

$ws.Range("H106").Value = 35900
$ws.Range("J106").Value = 35900
$ws.Range("L106").Value = 35900
$ws.Range("N106").Value = -38424

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H64").Value = 29999.666
$ws.Range("J64").Value = 29999.666
$ws.Range("L64").Value = 29999.666
$ws.Range("N64").Value = -30495.666

$ws.Range("H67").Value = 29999.666
$ws.Range("J67").Value = 29999.666
$ws.Range("L67").Value = 29999.666
$ws.Range("N67").Value = -31715.666

$ws.Range("H109").Value = 14950
$ws.Range("I109").Value = 14200
$ws.Range("J109").Value = 15200
$ws.Range("K109").Value = 14200
$ws.Range("L109").Value = 15200
$ws.Range("M109").Value = -12813
$ws.Range("N109").Value = -17974

$ws.Range("H115").Value = 44500
$ws.Range("J115").Value = 44500
$ws.Range("L115").Value = 44500
$ws.Range("N115").Value = -47634
